$d = $word.ActiveDocument

# --- locate the two paragraphs we need to touch, by content -------------
$metaIndex = 0
$imagePromptIndex = 0
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($metaIndex -eq 0 -and $t -like "Meta description:*") {
        $metaIndex = $i
    }
    if ($imagePromptIndex -eq 0 -and $t -like "Create a feature image*") {
        $imagePromptIndex = $i
    }
}

# 1. Remove the "Meta description: ..." paragraph that currently sits right
#    after the H1 title at the top of the document.
if ($metaIndex -gt 0) {
    $d.Paragraphs($metaIndex).Range.Delete()
}

# Re-resolve the image-prompt paragraph's index after the deletion above
# shifted everything that follows it up by one.
$count = $d.Paragraphs.Count
$imagePromptIndex = 0
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "Create a feature image*") {
        $imagePromptIndex = $i
    }
}

# 2. Insert a new bold paragraph carrying the title text
#    ("Play Coin-o-Mania Free Slot Game Review") right before the
#    "Create a feature image..." paragraph, matching the original
#    title/meta-description run layout (leading empty run + styled run).
$targetPara = $d.Paragraphs($imagePromptIndex)
$titleXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Coin-o-Mania Free Slot Game Review</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$targetPara.Range.InsertXML($titleXml)

# 3. Replace the text of the (now shifted) image-prompt paragraph with the
#    meta-description copy, keeping its leading empty run and italic
#    formatting on the text run intact.
$imagePromptIndex = $imagePromptIndex + 1
$finalPara = $d.Paragraphs($imagePromptIndex)
$descXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Find out about the pirate-themed slot game Coin-o-Mania, featuring bonus rounds and an RTP of 96%. Play for free now.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$finalRange = $d.Range($finalPara.Range.Start, $finalPara.Range.End - 1)
$finalRange.InsertXML($descXml)
